# Apply updated cryptocurrency price/volume data (commit: "Updated symbol list on Sat Jan 28 11:33:58 UTC 2023 with GitHub Actions")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell -> new text value (Price column D and Volume(1h) column E)
$updates = [ordered]@{
    "D2" = "307.24"
    "E2" = "0.64%"
    "D3" = "38.57"
    "E3" = "8.40%"
    "D4" = "5.098"
    "E4" = "1.03%"
    "D5" = "0.08112"
    "E5" = "1.34%"
    "D6" = "1.961"
    "E6" = "5.20%"
    "D7" = "7.949"
    "E7" = "2.13%"
    "D8" = "0.9294"
    "E8" = "0.83%"
    "D9" = "0.1416"
    "E9" = "10.15%"
    "D10" = "0.1954"
    "E10" = "2.73%"
    "D11" = "0.09012"
    "E11" = "-1.55%"
    "D12" = "0.03503"
    "E12" = "2.80%"
    "D13" = "0.09824"
    "E13" = "-0.39%"
    "D14" = "0.001406"
    "E14" = "0.30%"
    "D15" = "0.006044"
    "E15" = "-1.93%"
    "D16" = "3.731"
    "E16" = "-3.15%"
    "D17" = "4.197"
    "E17" = "1.52%"
    "E18" = "5.25%"
    "D19" = "0.3462"
    "E19" = "1.22%"
    "D20" = "0.1341"
    "E20" = "-0.50%"
    "D21" = "4.794"
    "E21" = "-7.73%"
    "D22" = "0.2454"
    "E22" = "6.28%"
    "D23" = "0.04404"
    "E23" = "-0.56%"
    "D24" = "0.001225"
    "E24" = "-0.73%"
    "E25" = "-1.05%"
    "E27" = "4.09%"
    "D39" = "0.02081"
    "E39" = "7.01%"
    "D40" = "0.05108"
    "E40" = "-1.47%"
    "D41" = "0.007481"
    "E41" = "-2.23%"
    "D42" = "0.01014"
    "E42" = "0.10%"
    "E43" = "0.44%"
    "D44" = "0.002134"
    "E44" = "-0.84%"
    "D45" = "0.009275"
    "E45" = "-3.66%"
    "D46" = "0.00006243"
    "E46" = "1.13%"
    "E47" = "0.07%"
    "D48" = "0.003029"
    "D49" = "0.001601"
    "E49" = "-3.59%"
    "E50" = "0.07%"
    "E51" = "0.07%"
}

foreach ($cellRef in $updates.Keys) {
    $cell = $ws.Range($cellRef)
    # Force text format so numeric-looking strings (e.g. "307.24") and
    # percent-looking strings (e.g. "0.64%") are stored as literal text,
    # matching the original inline-string cell type instead of being
    # auto-converted to a number/percentage by Excel.
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$cellRef]
    # Restore the default "General" number format so the cell format
    # matches the original (unformatted) cells.
    $cell.NumberFormat = "General"
}
